$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the species-record data among rows
# 7, 8, 10, 11, 12, 14 (columns A,B,D,E,F,G,H,Q,R). Rows 9 and 13 are
# untouched. Target values below are taken directly from the target
# (post-edit) state.

$rows = @{
    7  = @{ A = 112083125; B = 89369;  D = "LC"; E = 5447;   F = "Vedticka";        G = "Fuscoporia viticola";      H = "(Schwein.) Murrill";          Q = 413015.9403039298; R = 6656414.640994807 }
    8  = @{ A = 112083128; B = 77186;  D = "NT"; E = 353;    F = "Dvärgbägarlav";   G = "Cladonia parasitica";      H = "(Hoffm.) Hoffm.";             Q = 413190.1061828797; R = 6656475.01450387 }
    10 = @{ A = 112083110; B = 78107;  D = "NT"; E = 6453;   F = "Vedskivlav";      G = "Hertelidea botryosa";      H = "(Fr.) Printzen & Kantvilas";  Q = 412205.6393663768; R = 6656050.944565876 }
    11 = @{ A = 112083111; B = 90666;  D = "LC"; E = 4364;   F = "Dropptaggsvamp";  G = "Hydnellum ferrugineum";    H = "(Fr.:Fr.) P. Karst.";         Q = 412204.6634863199; R = 6655988.977203708 }
    12 = @{ A = 112083126; B = 78536;  D = "LC"; E = 229497; F = "Korallblylav";    G = "Parmeliella triptophylla"; H = "(Ach.) Müll.Arg.";            Q = 413016.7201701452; R = 6656341.641577623 }
    14 = @{ A = 112083118; B = 94134;  D = "NT"; E = 53;     F = "Vedtrappmossa";   G = "Crossocalyx hellerianus";  H = "(Nees ex Lindenb.) Meyl.";    Q = 412576.6879626553; R = 6656303.56951345 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
}
